$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) values are plain text in the source sheet (e.g. "42.089.79",
# "1.00"); a leading apostrophe forces Excel to store them as text instead of
# auto-converting to a number (which would drop formatting like trailing zeros).

$ws.Range("D2").Value = "'42.089.79"
$ws.Range("E2").Value = "  -1.97%  "
$ws.Range("D3").Value = "'2.263.13"
$ws.Range("E3").Value = "  -3.30%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'299.22"
$ws.Range("E5").Value = "  -2.43%  "
$ws.Range("D6").Value = "'94.13"
$ws.Range("E6").Value = "  -6.25%  "
$ws.Range("D7").Value = "'0.497"
$ws.Range("E7").Value = "  -3.02%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  -3.94%  "
$ws.Range("D10").Value = "'33.03"
$ws.Range("E10").Value = "  -5.62%  "
$ws.Range("D12").Value = "'47.89"
$ws.Range("E12").Value = "  -8.14%  "
$ws.Range("E13").Value = "  +0.11%  "
$ws.Range("E14").Value = "  -1.94%  "
$ws.Range("D15").Value = "'2.614.64"
$ws.Range("E15").Value = "  -3.34%  "
$ws.Range("D16").Value = "'15.43"
$ws.Range("E16").Value = "  -3.62%  "
$ws.Range("D17").Value = "'2.257.75"
$ws.Range("E17").Value = "  -4.24%  "
$ws.Range("E18").Value = "  -4.21%  "
$ws.Range("D19").Value = "'42.069.76"
$ws.Range("E19").Value = "  -1.84%  "
$ws.Range("E20").Value = "  -2.19%  "
$ws.Range("E21").Value = "  -3.37%  "
$ws.Range("D22").Value = "'11.34"
$ws.Range("E22").Value = "  -3.35%  "
$ws.Range("D23").Value = "'66.64"
$ws.Range("E23").Value = "  -1.97%  "
$ws.Range("D24").Value = "'233.52"
$ws.Range("E24").Value = "  -1.41%  "
$ws.Range("E25").Value = "  -4.87%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("E27").Value = "  -4.05%  "
$ws.Range("D28").Value = "'23.73"
$ws.Range("E28").Value = "  -7.50%  "
$ws.Range("D29").Value = "'2.16"
$ws.Range("E29").Value = "  -6.84%  "
$ws.Range("D30").Value = "'166.98"
$ws.Range("E30").Value = "  +3.36%  "
$ws.Range("D31").Value = "'33.69"
$ws.Range("E31").Value = "  -3.97%  "
$ws.Range("E32").Value = "  -3.16%  "
$ws.Range("D33").Value = "'1.00"
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("E34").Value = "  -3.86%  "
$ws.Range("E35").Value = "  -3.90%  "
$ws.Range("E36").Value = "  -5.74%  "
$ws.Range("D37").Value = "'0.0692"
$ws.Range("E37").Value = "  -5.00%  "
$ws.Range("B38").Value = "Celestia"
$ws.Range("C38").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D38").Value = "'16.14"
$ws.Range("E38").Value = "  -7.18%  "
$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").Value = "'2.78"
$ws.Range("E39").Value = "  -5.92%  "
$ws.Range("E40").Value = "  -3.55%  "
$ws.Range("E41").Value = "  -3.50%  "
$ws.Range("E42").Value = "  -8.56%  "
$ws.Range("E43").Value = "  -1.69%  "
$ws.Range("D44").Value = "'1.955.96"
$ws.Range("E44").Value = "  -2.88%  "
$ws.Range("E45").Value = "  -2.47%  "
$ws.Range("D46").Value = "'17.48"
$ws.Range("E46").Value = "  -6.96%  "
$ws.Range("D47").Value = "'9.53"
$ws.Range("E47").Value = "  -7.32%  "
$ws.Range("E48").Value = "  -4.67%  "
$ws.Range("E49").Value = "  -2.68%  "
$ws.Range("D50").Value = "'2.487.90"
$ws.Range("E50").Value = "  -2.69%  "
$ws.Range("D51").Value = "'51.71"
$ws.Range("E51").Value = "  -7.22%  "
